# Add new power plants to the Electricity Source subscript (issues #280 and #99)
#
# This adds six new technology rows to the "DCpUC" sheet, each pointing its
# Decommissioning Cost ($/MW) figure at an existing category already present
# on the sheet:
#   - hard coal w CCS                  -> same cost as "hard coal"      (B2)
#   - natural gas combined cycle w CCS -> same cost as "natural gas"    (B4)
#   - biomass w CCS                    -> same cost as "biomass"        (B10)
#   - lignite w CCS                    -> same cost as "lignite"        (B14)
#   - small modular reactor            -> same cost as "nuclear"        (B5)
#   - hydrogen                         -> same cost as "natural gas"    (B4)

$wb = $excel.ActiveWorkbook

$dcpuc = $wb.Worksheets.Item("DCpUC")
$data  = $wb.Worksheets.Item("Data")
$about = $wb.Worksheets.Item("About")

# New rows 19-24: label in column A, formula (copying an existing cost) in column B.
$newRows = @(
    @{ Row = 19; Label = "hard coal w CCS";                  Formula = "=B2"  },
    @{ Row = 20; Label = "natural gas combined cycle w CCS";  Formula = "=B4"  },
    @{ Row = 21; Label = "biomass w CCS";                     Formula = "=B10" },
    @{ Row = 22; Label = "lignite w CCS";                     Formula = "=B14" },
    @{ Row = 23; Label = "small modular reactor";             Formula = "=B5"  },
    @{ Row = 24; Label = "hydrogen";                          Formula = "=B4"  }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $dcpuc.Range("A$r").Value = $entry.Label
    $bCell = $dcpuc.Range("B$r")
    $bCell.Formula = $entry.Formula
    # Match the integer ("0") number format used by the rest of column B.
    $bCell.NumberFormat = "0"
}

# Incidental cleanup: the stray (emptied, but still formatted) cell on the
# Data sheet is fully cleared so it no longer occupies a cell record.
$data.Range("C15").Clear() | Out-Null

# Restore view state: DCpUC keeps its own selection memory at A25 (just past
# the newly-added rows) while "About" remains the active sheet/tab with its
# selection at B40.
$dcpuc.Range("A25").Select() | Out-Null
$about.Activate() | Out-Null
$about.Range("B40").Select() | Out-Null
